$wb = $excel.ActiveWorkbook

# --- start_price sheet: A2 ---
$wsStartPrice = $wb.Worksheets.Item("start_price")
$wsStartPrice.Range("A2").Value = 2.404

# --- Linear sheet: mu (B2), B (B3), sig2 (B4), abs_epsi_autocorr (B5) ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.0004149177480496889
$wsLinear.Range("B3").Value = -0.01905430008293248
$wsLinear.Range("B4").Value = 0.02966496776916974
$wsLinear.Range("B5").Value = "[1.0, 0.2158584620117635, 0.07921189633648493, 0.08008533595227837, 0.04813017056828746, 0.06362111040704166, 0.21889671202069694, 0.3733268026268109, 0.2155397340088252, 0.06434753572771686, 0.02538232397625779, 0.04064795515769978, 0.06957247562040036, 0.21050128680010893, 0.36416141511695754, 0.22413347954352952, 0.028313399671818423, 0.03939200844539863, 0.03569433084334304, 0.049998537702642]"

# --- NonLinear sheet: p (B3), mu_0 (B4), B_0 (B5), sig2_0 (B6), mu_1 (B7), B_1 (B8), sig2_1 (B9) ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 1.060405156537753
$wsNonLinear.Range("B4").Value = 0.007197046582406103
$wsNonLinear.Range("B5").Value = 0.1121724628569274
$wsNonLinear.Range("B6").Value = 0.02549378743211238
$wsNonLinear.Range("B7").Value = 0.002045731721457663
$wsNonLinear.Range("B8").Value = -0.09420429781706366
$wsNonLinear.Range("B9").Value = 0.03403584203939997
$wsNonLinear.Range("B10").Value = "[1.0, 0.21490630963217355, 0.08582636982261824, 0.08983001010006395, 0.057518836548398344, 0.0715506465272342, 0.22038742740653314, 0.3685917942213751, 0.2144502988766464, 0.07174527256239713, 0.03272487559404901, 0.05006359524982565, 0.07674398819209349, 0.21078048932175994, 0.36056536290154806, 0.22553605174698974, 0.03558746497814142, 0.04621397598906925, 0.04289470404750156, 0.055505233613785974]"
